$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Formula = "=LOG10(17.2)"
$ws.Range("C2").Formula = "=LOG10(0.355/3.07)"
$ws.Range("D2").Formula = "=LOG10(0.0355/3.07)"
$ws.Range("E2").Formula = "=LOG10(1.33/3.07)"

# Row 3
$ws.Range("B3").Formula = "=LOG10(13.9)"
$ws.Range("C3").Formula = "=LOG10(0.342/3.07)"
$ws.Range("D3").Formula = "=LOG10(0.045/3.07)"
$ws.Range("E3").Formula = "=LOG10(1.22/3.07)"

# Row 4 (B4 unchanged)
$ws.Range("C4").Formula = "=LOG10(0.401/3.07)"
$ws.Range("D4").Formula = "=LOG10(0.0484/3.07)"
$ws.Range("E4").Formula = "=LOG10(1.28/3.07)"

# Row 5
$ws.Range("B5").Formula = "=LOG10(24.3)"
$ws.Range("C5").Formula = "=LOG10(0.241/3.09)"
$ws.Range("D5").Formula = "=LOG10(0.0798/3.09)"
$ws.Range("E5").Formula = "=LOG10(1.76/3.09)"

# Update selection to D2
$ws.Range("D2").Select()
